$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 38
$ws.Range("E10").Value = 344
$ws.Range("E11").Value = 238
$ws.Range("E13").Value = 101
$ws.Range("E22").Value = 133
$ws.Range("E23").Value = 156
$ws.Range("E26").Value = 97
$ws.Range("E27").Value = 238
$ws.Range("E30").Value = 158
$ws.Range("E31").Value = 62
$ws.Range("E34").Value = 161
$ws.Range("E39").Value = 151
$ws.Range("E42").Value = 257
$ws.Range("E44").Value = 240
$ws.Range("E47").Value = 325
$ws.Range("E50").Value = 185
$ws.Range("F50").Value = 64
$ws.Range("H50").Value = 64
